# "IP and Licences as separate episode"
#
# 1. Delete the "License, Copyright and Data" slide (it is being split out
#    into its own, separate episode) so the deck goes from 19 -> 18 slides.
# 2. Refresh the cached "datetimeFigureOut" footer date (14/10/2021 ->
#    18/10/2021) everywhere it is cached: the slide master and every
#    slide layout.

$p = $ppt.ActivePresentation

# --- 1. Remove the "License, Copyright and Data" slide -------------------
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    $slideTitle = ""
    foreach ($sh in $s.Shapes) {
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $slideTitle = $sh.TextFrame.TextRange.Text
                break
            }
        }
    }
    if ($slideTitle -like "License*") {
        $s.Delete()
    }
}

# --- 2. Update the cached footer date on the master + every layout -------
$newDate = "18/10/2021"

$master = $p.SlideMaster
foreach ($sh in $master.Shapes) {
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "14/10/2021") {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    foreach ($sh in $layout.Shapes) {
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq "14/10/2021") {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}
